$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers (human-readable labels)
$ws.Range("B1").Value = "Depuradoras"
$ws.Range("C1").Value = "Habitantes equivalentes"
$ws.Range("D1").Value = "Comarca nombre"
$ws.Range("E1").Value = "Comarca código"
$ws.Range("F1").Value = "CCAA código"
$ws.Range("G1").Value = "Municipio código"
$ws.Range("H1").Value = "Año"
$ws.Range("I1").Value = "Municipio nombre"

# Row 2 - concept / dimension / measure identifiers
$ws.Range("A2").Value = "iaest-measure:"
$ws.Range("B2").Value = "iaest-measure:depuradoras"
$ws.Range("C2").Value = "iaest-measure:habitantes-equivalentes"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "iaest-measure:habitantes-equivalentes"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "sdmx-dimension:refPeriod"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3 - role (medida/dim)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "dim"
$ws.Range("I3").Value = "dim"

# Row 4 - data type / codelist
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-comarca"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "xsd:date"
$ws.Range("I4").Value = "URI-Municipio"

# Row 5 - B5 no longer used; H5 carries the mapping file reference now
$ws.Range("B5").ClearContents()
$ws.Range("H5").Value = "mapping-ano.xlsx"
